# Rename the sheet from RentBikeServiceController -> Bike
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Bike"

# Update the "Method name:" header in row 3 for the Bike constructor section
$ws.Range("A3").Value = "Method name: Bike (Constructor)"
$ws.Range("A3").Font.Bold = $true

# Row 5: parameter "name"
$ws.Range("A5").Value = 1
$ws.Range("B5").Value = "name"
$ws.Range("C5").Value = "'- must not be null or empty`n- must start with a letter`n- can only contain letters, number, space, dash - and underscore _"
$ws.Range("C5").WrapText = $true
$ws.Rows(5).RowHeight = 71.25

# Row 6: parameter "bike type" (bold row, like the original header formatting)
$ws.Range("B6").Value = "bike type"
$ws.Range("B6").Font.Bold = $true
$ws.Range("C6").Value = "'- must not be null or empty`n- must start with a letter`n- can only contain letters, number, space, dash - and underscore _"
$ws.Range("C6").WrapText = $true
$ws.Rows(6).RowHeight = 72

# Row 7: parameter "bike image"
$ws.Range("B7").Value = "bike image"
$ws.Range("B7").Font.Bold = $true

# Row 8: parameter "bike barcode"
$ws.Range("B8").Value = "bike barcode"
$ws.Range("B8").Font.Bold = $true

# Row 9: parameter "bike rental price"
$ws.Range("B9").Value = "bike rental price"
$ws.Range("B9").Font.Bold = $true
$ws.Range("C9").Value = "'- must greater than 0"

# Row 10: parameter "deposit price"
$ws.Range("B10").Value = "deposit price"
$ws.Range("B10").Font.Bold = $true
$ws.Range("C10").Value = "'- must greater than 0"

# Row 11: parameter "currency"
$ws.Range("B11").Value = "currency"
$ws.Range("B11").Font.Bold = $true
$ws.Range("C11").Value = "'- must of 3 letter, abbreviation"

# Row 12: parameter "create date"
$ws.Range("B12").Value = "create date"
$ws.Range("B12").Font.Bold = $true

# Row 13: parameter "total rental time"
$ws.Range("B13").Value = "total rental time"
$ws.Range("B13").Font.Bold = $true
$ws.Range("C13").Value = "'- must greater than 0"

# Row 14: parameter "status"
$ws.Range("B14").Value = "status"
$ws.Range("B14").Font.Bold = $true
$ws.Range("C14").Value = "'- can only be FREE or RENTED"

# Row 15: parameter "battery"
$ws.Range("B15").Value = "battery"
$ws.Range("B15").Font.Bold = $true
$ws.Range("C15").Value = "'- must be in range 0 - 100"

# Row 18/19: next method header block (Dock), moved down from rows 6/7
$ws.Range("A18").Value = "Method name:"
$ws.Range("A18").Font.Bold = $true
$ws.Range("B18").Value = ""
$ws.Range("B18").Font.Bold = $true
$ws.Range("C18").Value = ""
$ws.Range("C18").Font.Bold = $true

$ws.Range("A19").Value = "#"
$ws.Range("B19").Value = "ParamName"
$ws.Range("C19").Value = "Conditions"
$ws.Range("A19:C19").Font.Bold = $true
$ws.Range("A19:C19").Borders.LineStyle = 1

$ws.Range("C16").Select()
